$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")
$ws.Activate()

# Data for the new rows (491-501), matching the existing column layout:
# A=DATE, B=VECHILE REG NO, C=VEHICLE BRAND, D=ISSUE, E=STATUS, F=AMOUNT, G=CASH TYPE
$rows = @(
    @{ Row=491; Date=44810; B="KA03MV7501"; C="VENTO";    D="PMS";                 E="WORK DONE DELIVERED"; F=2998; G="CREDIT" },
    @{ Row=492; Date=44810; B="UP14CK8538"; C="SAEL";     D="RR";                  E="WORK DONE DELIVERED"; F=3250; G=$null },
    @{ Row=493; Date=44810; B="PB35R0757";  C="SX4";      D="BODY SHOP";           E="WORK IN PROGRESS";    F=$null; G=$null },
    @{ Row=494; Date=44810; B="KA53Z5764";  C="LIEA";     D="BODY SHOP";           E="WORK IN PROGRESS";    F=$null; G=$null },
    @{ Row=495; Date=44811; B="KA53N0153";  C="SPARK";    D="WIPER BLADE CHANGE";  E="WORK DONE DELIVERED"; F=500;  G="CASH" },
    @{ Row=496; Date=44811; B="TN22CM8979"; C="MICRA";    D="RUNNING REPAIR";      E="WORK IN PROGRESS";    F=$null; G=$null },
    @{ Row=497; Date=44811; B="KA03MX2103"; C="WAGON R";  D="PMS";                 E="WORK IN PROGRESS";    F=$null; G=$null },
    @{ Row=498; Date=44811; B="KA53MJ8304"; C="I20";      D="RAT MESH";            E="WORK DONE DELIVERED"; F=4000; G="PAYTM" },
    @{ Row=499; Date=44811; B="KA51MB3247"; C="POLO";     D="SCANNING";            E="WORK DONE DELIVERED"; F=1000; G="PAYTM" },
    @{ Row=500; Date=44811; B="KA08M3568";  C="ALTO 800"; D="SILENCER WELDING";    E="WORK DONE DELIVERED"; F=800;  G="PAYTM" },
    @{ Row=501; Date=44811; B="KA01MM2572"; C="INNOVA";   D="AC LEAKAGE PROBLEM "; E="WORK IN PROGRESS";    F=$null; G=$null }
)

# Reference cells (row 490 is the last existing data row) whose number
# formatting we reuse for the new date / amount cells, so no new numFmt /
# cellXfs entries get invented.
$dateRefCell = $ws.Cells.Item(490, 1)
$amountRefCell = $ws.Cells.Item(490, 6)

foreach ($r in $rows) {
    $rowNum = $r.Row

    $dateRefCell.Copy()
    $ws.Cells.Item($rowNum, 1).PasteSpecial(-4122)
    $ws.Cells.Item($rowNum, 1).Value = $r.Date

    $ws.Cells.Item($rowNum, 2).Value = $r.B
    $ws.Cells.Item($rowNum, 3).Value = $r.C
    $ws.Cells.Item($rowNum, 4).Value = $r.D
    $ws.Cells.Item($rowNum, 5).Value = $r.E

    if ($null -ne $r.F) {
        $amountRefCell.Copy()
        $ws.Cells.Item($rowNum, 6).PasteSpecial(-4122)
        $ws.Cells.Item($rowNum, 6).Value = $r.F
    }
    if ($null -ne $r.G) {
        $ws.Cells.Item($rowNum, 7).Value = $r.G
    }
}

$excel.CutCopyMode = 0

$ws.Range("F501").Select()
$excel.ActiveWindow.ScrollRow = 487
